# Add data-expand result:
#  - E11:J11 currently hold numeric 0's; replace them with the text "zeros"
#    (new shared string), styled with the (new) 11pt 宋体 font.
#  - Move the sheet's active-cell selection from D17 to E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("E11:J11")
$range.Value = "zeros"
$range.Font.Name = "宋体"

$ws.Range("E13").Select()
